# Add two new columns "I0" (column I) and "IF" (column J) to the sheet.
# The new header cells should share the same formatting (bold, centered,
# bordered) already used by the other header cells in row 1, so we copy
# the formatting from the existing H1 header cell rather than building a
# brand-new style.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers (row 1)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false

# Data values for rows 2-14 (columns I = I0, J = IF)
$data = @{
    2  = @(6, 7)
    3  = @(8, 8)
    4  = @(6, 6)
    5  = @(6, 7)
    6  = @(7, 7)
    7  = @(7, 8)
    8  = @(13, 13)
    9  = @(9, 9)
    10 = @(6, 7)
    11 = @(7, 8)
    12 = @(7, 7)
    13 = @(8, 8)
    14 = @(9, 9)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 9).Value = $vals[0]
    $ws.Cells.Item($row, 10).Value = $vals[1]
}
